{"js": "const doc = context.document;\nconst body = doc.body;\n\n// The \"_GoBack\" bookmark currently sits right before the trailing\n// \" - master branch\" run of the (only) paragraph. It needs to move to\n// the new last paragraph, so drop it from its current spot first.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Grab the paragraph that holds \"File02 - master branch\" so we can anchor\n// the two new paragraphs after it.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst firstPara = body.paragraphs.items[body.paragraphs.items.length - 1];\n\n// Insert a blank paragraph right after it...\nconst blankPara = firstPara.insertParagraph(\"\", Word.InsertLocation.after);\n\n// ...followed by a paragraph with the new sentence.\nconst newPara = blankPara.insertParagraph(\n  \"Make second change to file02 in master branch\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Re-query the paragraphs so the range we grab next reflects the content\n// that was just inserted (a stale Paragraph object yields a range spanning\n// the whole paragraph instead of a collapsed point).\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst refreshedItems = body.paragraphs.items;\nconst refreshedNewPara = refreshedItems[refreshedItems.length - 1];\n\n// Re-insert the bookmark, collapsed right at the start of the new\n// last paragraph (before its text), matching its original placement style.\nconst startRange = refreshedNewPara.getRange(Word.RangeLocation.start);\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently sits right before the trailing\n# \" - master branch\" run of the (only) paragraph. It needs to move to the\n# new last paragraph, so remove it from its current spot first.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Grab the paragraph that holds \"File02 - master branch\" so we can anchor\n# the two new paragraphs after it.\n$firstPara = $d.Paragraphs.Item(1)\n$endOfFirst = $firstPara.Range\n$endOfFirst.Collapse(0)  # wdCollapseEnd\n\n# Insert a blank paragraph right after it ...\n$endOfFirst.InsertParagraphAfter()\n\n# ... then split off a further paragraph right after that blank one, and\n# fill it with the new sentence (inserting straight into the blank\n# paragraph would just append the text to it instead of creating a new one).\n$secondPara = $d.Paragraphs.Item(2)\n$endOfSecond = $secondPara.Range\n$endOfSecond.Collapse(0)  # wdCollapseEnd\n$endOfSecond.InsertParagraphAfter()\n\n$thirdPara = $d.Paragraphs.Item(3)\n$thirdPara.Range.InsertAfter(\"Make second change to file02 in master branch\")\n\n# Re-insert the bookmark, collapsed right at the start of the new last\n# paragraph (before its text), matching its original placement style.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$startOfLast = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)\n$d.Bookmarks.Add(\"_GoBack\", $startOfLast) | Out-Null\n"}
